$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata ---
$ws.Range("A1:AJ189").Select() | Out-Null

# --- Cell value updates (CRI-E_T / CRI-W_T / CRI-F_T fill-in + CRI-Total corrections) ---
$ws.Cells.Item(3, 28).Value = 40
$ws.Cells.Item(3, 32).Value = 113
$ws.Cells.Item(5, 32).Value = 99
$ws.Cells.Item(7, 32).Value = 126
$ws.Cells.Item(12, 32).Value = 124
$ws.Cells.Item(14, 32).Value = 110
$ws.Cells.Item(15, 32).Value = 124
$ws.Cells.Item(16, 32).Value = 121
$ws.Cells.Item(17, 28).Value = 40
$ws.Cells.Item(17, 32).Value = 138
$ws.Cells.Item(18, 32).Value = 123
$ws.Cells.Item(19, 28).Value = 50
$ws.Cells.Item(19, 32).Value = 159
$ws.Cells.Item(21, 32).Value = 135
$ws.Cells.Item(22, 28).Value = 20
$ws.Cells.Item(23, 32).Value = 135
$ws.Cells.Item(26, 32).Value = 123
$ws.Cells.Item(27, 28).Value = 40
$ws.Cells.Item(27, 32).Value = 133
$ws.Cells.Item(31, 28).Value = 40
$ws.Cells.Item(31, 32).Value = 153
$ws.Cells.Item(32, 32).Value = 142
$ws.Cells.Item(35, 28).Value = 30
$ws.Cells.Item(41, 28).Value = 40
$ws.Cells.Item(41, 32).Value = 149
$ws.Cells.Item(43, 28).Value = 115
$ws.Cells.Item(43, 32).Value = 149
$ws.Cells.Item(44, 27).Value = 20
$ws.Cells.Item(44, 28).Value = 15
$ws.Cells.Item(45, 27).Value = 15
$ws.Cells.Item(45, 28).Value = 15
$ws.Cells.Item(45, 32).Value = 104
$ws.Cells.Item(46, 27).Value = 20
$ws.Cells.Item(46, 28).Value = 10
$ws.Cells.Item(47, 27).Value = 20
$ws.Cells.Item(47, 28).Value = 10
$ws.Cells.Item(48, 27).Value = 20
$ws.Cells.Item(48, 28).Value = 5
$ws.Cells.Item(49, 27).Value = 20
$ws.Cells.Item(49, 28).Value = 50
$ws.Cells.Item(49, 32).Value = 136
$ws.Cells.Item(50, 27).Value = 20
$ws.Cells.Item(50, 28).Value = 25
$ws.Cells.Item(51, 27).Value = 25
$ws.Cells.Item(51, 28).Value = 35
$ws.Cells.Item(51, 32).Value = 131
$ws.Cells.Item(52, 27).Value = 15
$ws.Cells.Item(52, 28).Value = 40
$ws.Cells.Item(52, 32).Value = 136
$ws.Cells.Item(53, 27).Value = 21
$ws.Cells.Item(53, 28).Value = 0
$ws.Cells.Item(54, 27).Value = 15
$ws.Cells.Item(54, 28).Value = 55
$ws.Cells.Item(54, 32).Value = 138
$ws.Cells.Item(55, 27).Value = 20
$ws.Cells.Item(55, 28).Value = 45
$ws.Cells.Item(55, 32).Value = 138
$ws.Cells.Item(56, 27).Value = 26
$ws.Cells.Item(56, 28).Value = 5
$ws.Cells.Item(57, 27).Value = 24
$ws.Cells.Item(57, 28).Value = 40
$ws.Cells.Item(57, 32).Value = 125
$ws.Cells.Item(58, 27).Value = 2
$ws.Cells.Item(58, 28).Value = 35
$ws.Cells.Item(59, 27).Value = 20
$ws.Cells.Item(59, 28).Value = 10
$ws.Cells.Item(60, 27).Value = 26
$ws.Cells.Item(60, 28).Value = 10
$ws.Cells.Item(61, 27).Value = 20
$ws.Cells.Item(61, 28).Value = 35
$ws.Cells.Item(61, 32).Value = 136
$ws.Cells.Item(62, 27).Value = 20
$ws.Cells.Item(62, 28).Value = 50
$ws.Cells.Item(63, 27).Value = 20
$ws.Cells.Item(63, 28).Value = 50
$ws.Cells.Item(64, 27).Value = 20
$ws.Cells.Item(64, 28).Value = 40
$ws.Cells.Item(65, 27).Value = 20
$ws.Cells.Item(65, 28).Value = 30
$ws.Cells.Item(66, 27).Value = 20
$ws.Cells.Item(66, 28).Value = 25
$ws.Cells.Item(67, 27).Value = 20
$ws.Cells.Item(67, 28).Value = 15
$ws.Cells.Item(68, 27).Value = 15
$ws.Cells.Item(68, 28).Value = 30
$ws.Cells.Item(69, 27).Value = 18
$ws.Cells.Item(69, 28).Value = 35
$ws.Cells.Item(69, 32).Value = 133
$ws.Cells.Item(70, 27).Value = 21
$ws.Cells.Item(70, 28).Value = 40
$ws.Cells.Item(70, 32).Value = 110
$ws.Cells.Item(71, 27).Value = 23
$ws.Cells.Item(71, 28).Value = 35
$ws.Cells.Item(71, 32).Value = 129
$ws.Cells.Item(72, 27).Value = 20
$ws.Cells.Item(72, 28).Value = 45
$ws.Cells.Item(72, 32).Value = 113
$ws.Cells.Item(73, 27).Value = 26
$ws.Cells.Item(73, 28).Value = 110
$ws.Cells.Item(73, 32).Value = 151
$ws.Cells.Item(74, 27).Value = 30
$ws.Cells.Item(74, 28).Value = 70
$ws.Cells.Item(74, 32).Value = 166
$ws.Cells.Item(75, 27).Value = 20
$ws.Cells.Item(75, 28).Value = 35
$ws.Cells.Item(75, 32).Value = 130
$ws.Cells.Item(76, 27).Value = 20
$ws.Cells.Item(76, 28).Value = 65
$ws.Cells.Item(76, 32).Value = 136
$ws.Cells.Item(77, 27).Value = 27
$ws.Cells.Item(77, 28).Value = 50
$ws.Cells.Item(77, 32).Value = 148
$ws.Cells.Item(78, 27).Value = 29
$ws.Cells.Item(78, 28).Value = 45
$ws.Cells.Item(78, 29).Value = 257
$ws.Cells.Item(78, 32).Value = 144
$ws.Cells.Item(79, 27).Value = 21
$ws.Cells.Item(79, 28).Value = 15
$ws.Cells.Item(79, 29).Value = 40
$ws.Cells.Item(80, 27).Value = 21
$ws.Cells.Item(80, 28).Value = 35
$ws.Cells.Item(80, 29).Value = 332
$ws.Cells.Item(80, 32).Value = 132
$ws.Cells.Item(81, 27).Value = 29
$ws.Cells.Item(81, 28).Value = 50
$ws.Cells.Item(81, 29).Value = 270
$ws.Cells.Item(82, 27).Value = 20
$ws.Cells.Item(82, 28).Value = 10
$ws.Cells.Item(82, 29).Value = 40
$ws.Cells.Item(83, 27).Value = 20
$ws.Cells.Item(83, 28).Value = 5
$ws.Cells.Item(83, 29).Value = 35
$ws.Cells.Item(84, 27).Value = 20
$ws.Cells.Item(84, 28).Value = 25
$ws.Cells.Item(84, 29).Value = 232
$ws.Cells.Item(84, 32).Value = 117
$ws.Cells.Item(85, 27).Value = 15
$ws.Cells.Item(85, 28).Value = 35
$ws.Cells.Item(85, 29).Value = 242
$ws.Cells.Item(85, 32).Value = 119
$ws.Cells.Item(86, 27).Value = 20
$ws.Cells.Item(86, 28).Value = 40
$ws.Cells.Item(86, 29).Value = 248
$ws.Cells.Item(86, 32).Value = 125
$ws.Cells.Item(87, 27).Value = 25
$ws.Cells.Item(87, 28).Value = 40
$ws.Cells.Item(87, 29).Value = 278
$ws.Cells.Item(87, 32).Value = 133
$ws.Cells.Item(88, 27).Value = 15
$ws.Cells.Item(88, 28).Value = 45
$ws.Cells.Item(88, 29).Value = 363
$ws.Cells.Item(88, 32).Value = 129
$ws.Cells.Item(89, 27).Value = 20
$ws.Cells.Item(89, 28).Value = 10
$ws.Cells.Item(89, 29).Value = 35
$ws.Cells.Item(90, 27).Value = 20
$ws.Cells.Item(90, 28).Value = 15
$ws.Cells.Item(90, 29).Value = 25
$ws.Cells.Item(91, 27).Value = 15
$ws.Cells.Item(91, 28).Value = 10
$ws.Cells.Item(91, 29).Value = 55
$ws.Cells.Item(92, 27).Value = 15
$ws.Cells.Item(92, 28).Value = 5
$ws.Cells.Item(92, 29).Value = 80
$ws.Cells.Item(93, 27).Value = 30
$ws.Cells.Item(93, 28).Value = 50
$ws.Cells.Item(93, 29).Value = 367
$ws.Cells.Item(93, 32).Value = 158
$ws.Cells.Item(94, 27).Value = 20
$ws.Cells.Item(94, 28).Value = 60
$ws.Cells.Item(94, 29).Value = 289
$ws.Cells.Item(94, 32).Value = 135
$ws.Cells.Item(95, 27).Value = 20
$ws.Cells.Item(95, 28).Value = 15
$ws.Cells.Item(95, 29).Value = 80
$ws.Cells.Item(96, 27).Value = 20
$ws.Cells.Item(96, 28).Value = 5
$ws.Cells.Item(96, 29).Value = 35
$ws.Cells.Item(97, 27).Value = 35
$ws.Cells.Item(97, 28).Value = 40
$ws.Cells.Item(97, 29).Value = 257
$ws.Cells.Item(97, 32).Value = 150
$ws.Cells.Item(98, 27).Value = 20
$ws.Cells.Item(98, 28).Value = 25
$ws.Cells.Item(98, 29).Value = 130
$ws.Cells.Item(99, 27).Value = 20
$ws.Cells.Item(99, 28).Value = 5
$ws.Cells.Item(99, 29).Value = 60
$ws.Cells.Item(100, 27).Value = 15
$ws.Cells.Item(100, 28).Value = 5
$ws.Cells.Item(100, 29).Value = 25
$ws.Cells.Item(101, 27).Value = 15
$ws.Cells.Item(101, 28).Value = 5
$ws.Cells.Item(101, 29).Value = 25
$ws.Cells.Item(102, 27).Value = 15
$ws.Cells.Item(102, 28).Value = 5
$ws.Cells.Item(102, 29).Value = 35
$ws.Cells.Item(103, 27).Value = 15
$ws.Cells.Item(103, 28).Value = 35
$ws.Cells.Item(103, 29).Value = 87
$ws.Cells.Item(103, 32).Value = 100
$ws.Cells.Item(104, 27).Value = 20
$ws.Cells.Item(104, 28).Value = 10
$ws.Cells.Item(104, 29).Value = 40
$ws.Cells.Item(105, 27).Value = 20
$ws.Cells.Item(105, 28).Value = 5
$ws.Cells.Item(105, 29).Value = 35
$ws.Cells.Item(106, 27).Value = 15
$ws.Cells.Item(106, 28).Value = 5
$ws.Cells.Item(106, 29).Value = 35
$ws.Cells.Item(107, 27).Value = 20
$ws.Cells.Item(107, 28).Value = 5
$ws.Cells.Item(107, 29).Value = 26
$ws.Cells.Item(107, 32).Value = 98
$ws.Cells.Item(108, 27).Value = 15
$ws.Cells.Item(108, 28).Value = 5
$ws.Cells.Item(108, 29).Value = 33
$ws.Cells.Item(109, 27).Value = 30
$ws.Cells.Item(109, 28).Value = 15
$ws.Cells.Item(109, 29).Value = 105
$ws.Cells.Item(110, 27).Value = 30
$ws.Cells.Item(110, 28).Value = 20
$ws.Cells.Item(110, 29).Value = 150
$ws.Cells.Item(111, 27).Value = 20
$ws.Cells.Item(111, 28).Value = 5
$ws.Cells.Item(111, 29).Value = 55
$ws.Cells.Item(112, 27).Value = 20
$ws.Cells.Item(112, 28).Value = 35
$ws.Cells.Item(112, 29).Value = 247
$ws.Cells.Item(112, 32).Value = 119
$ws.Cells.Item(113, 27).Value = 40
$ws.Cells.Item(113, 28).Value = 35
$ws.Cells.Item(113, 29).Value = 463
$ws.Cells.Item(113, 32).Value = 165
$ws.Cells.Item(114, 27).Value = 27
$ws.Cells.Item(114, 28).Value = 30
$ws.Cells.Item(114, 29).Value = 307
$ws.Cells.Item(114, 32).Value = 128
$ws.Cells.Item(115, 27).Value = 29
$ws.Cells.Item(115, 28).Value = 40
$ws.Cells.Item(115, 29).Value = 233
$ws.Cells.Item(116, 27).Value = 20
$ws.Cells.Item(116, 28).Value = 5
$ws.Cells.Item(116, 29).Value = 90
$ws.Cells.Item(117, 27).Value = 26
$ws.Cells.Item(117, 28).Value = 0
$ws.Cells.Item(117, 29).Value = 85
$ws.Cells.Item(118, 27).Value = 22
$ws.Cells.Item(118, 28).Value = 30
$ws.Cells.Item(118, 29).Value = 648
$ws.Cells.Item(118, 32).Value = 150
$ws.Cells.Item(119, 27).Value = 18
$ws.Cells.Item(119, 28).Value = 40
$ws.Cells.Item(119, 29).Value = 215
$ws.Cells.Item(120, 27).Value = 20
$ws.Cells.Item(120, 28).Value = 45
$ws.Cells.Item(120, 29).Value = 410
$ws.Cells.Item(121, 27).Value = 20
$ws.Cells.Item(121, 28).Value = 10
$ws.Cells.Item(121, 29).Value = 80
$ws.Cells.Item(122, 27).Value = 20
$ws.Cells.Item(122, 28).Value = 5
$ws.Cells.Item(122, 29).Value = 55
$ws.Cells.Item(123, 27).Value = 15
$ws.Cells.Item(123, 28).Value = 5
$ws.Cells.Item(123, 29).Value = 30
$ws.Cells.Item(124, 27).Value = 20
$ws.Cells.Item(124, 28).Value = 5
$ws.Cells.Item(124, 29).Value = 40
$ws.Cells.Item(125, 27).Value = 20
$ws.Cells.Item(125, 28).Value = 40
$ws.Cells.Item(125, 29).Value = 313
$ws.Cells.Item(125, 32).Value = 134
$ws.Cells.Item(126, 27).Value = 10
$ws.Cells.Item(126, 28).Value = 40
$ws.Cells.Item(126, 29).Value = 298
$ws.Cells.Item(126, 32).Value = 118
$ws.Cells.Item(127, 27).Value = 20
$ws.Cells.Item(127, 28).Value = 25
$ws.Cells.Item(127, 29).Value = 252
$ws.Cells.Item(127, 32).Value = 123
$ws.Cells.Item(128, 27).Value = 43
$ws.Cells.Item(128, 28).Value = 55
$ws.Cells.Item(128, 29).Value = 471
$ws.Cells.Item(128, 32).Value = 177
$ws.Cells.Item(129, 27).Value = 17
$ws.Cells.Item(129, 28).Value = 40
$ws.Cells.Item(129, 29).Value = 302
$ws.Cells.Item(129, 32).Value = 126
$ws.Cells.Item(130, 27).Value = 25
$ws.Cells.Item(130, 28).Value = 25
$ws.Cells.Item(130, 29).Value = 125
$ws.Cells.Item(131, 27).Value = 20
$ws.Cells.Item(131, 28).Value = 40
$ws.Cells.Item(131, 29).Value = 400
$ws.Cells.Item(131, 32).Value = 138
$ws.Cells.Item(132, 27).Value = 12
$ws.Cells.Item(132, 28).Value = 40
$ws.Cells.Item(132, 29).Value = 223
$ws.Cells.Item(132, 32).Value = 107
$ws.Cells.Item(133, 27).Value = 18
$ws.Cells.Item(133, 28).Value = 30
$ws.Cells.Item(133, 29).Value = 202
$ws.Cells.Item(133, 32).Value = 114
$ws.Cells.Item(134, 27).Value = 20
$ws.Cells.Item(134, 28).Value = 10
$ws.Cells.Item(134, 29).Value = 40
$ws.Cells.Item(135, 27).Value = 22
$ws.Cells.Item(135, 28).Value = 35
$ws.Cells.Item(135, 29).Value = 287
$ws.Cells.Item(135, 32).Value = 138
$ws.Cells.Item(136, 27).Value = 23
$ws.Cells.Item(136, 28).Value = 30
$ws.Cells.Item(136, 29).Value = 251
$ws.Cells.Item(137, 27).Value = 20
$ws.Cells.Item(137, 28).Value = 5
$ws.Cells.Item(137, 29).Value = 40
$ws.Cells.Item(138, 27).Value = 15
$ws.Cells.Item(138, 28).Value = 5
$ws.Cells.Item(138, 29).Value = 40
$ws.Cells.Item(139, 27).Value = 22
$ws.Cells.Item(139, 28).Value = 45
$ws.Cells.Item(139, 29).Value = 353
$ws.Cells.Item(139, 32).Value = 132
$ws.Cells.Item(140, 27).Value = 20
$ws.Cells.Item(140, 28).Value = 50
$ws.Cells.Item(140, 29).Value = 408
$ws.Cells.Item(140, 32).Value = 146
$ws.Cells.Item(141, 27).Value = 25
$ws.Cells.Item(141, 28).Value = 30
$ws.Cells.Item(141, 29).Value = 377
$ws.Cells.Item(141, 32).Value = 134
$ws.Cells.Item(142, 27).Value = 20
$ws.Cells.Item(142, 28).Value = 35
$ws.Cells.Item(142, 29).Value = 384
$ws.Cells.Item(142, 32).Value = 125
$ws.Cells.Item(143, 27).Value = 25
$ws.Cells.Item(143, 28).Value = 35
$ws.Cells.Item(143, 29).Value = 188
$ws.Cells.Item(143, 32).Value = 124
$ws.Cells.Item(144, 27).Value = 15
$ws.Cells.Item(144, 28).Value = 15
$ws.Cells.Item(144, 29).Value = 142
$ws.Cells.Item(144, 32).Value = 100
$ws.Cells.Item(145, 27).Value = 30
$ws.Cells.Item(145, 28).Value = 35
$ws.Cells.Item(145, 29).Value = 178
$ws.Cells.Item(146, 27).Value = 31
$ws.Cells.Item(146, 28).Value = 5
$ws.Cells.Item(146, 29).Value = 416
$ws.Cells.Item(147, 27).Value = 25
$ws.Cells.Item(147, 28).Value = 20
$ws.Cells.Item(147, 29).Value = 95
$ws.Cells.Item(148, 27).Value = 31
$ws.Cells.Item(148, 28).Value = 35
$ws.Cells.Item(148, 29).Value = 272
$ws.Cells.Item(148, 32).Value = 131
$ws.Cells.Item(149, 27).Value = 30
$ws.Cells.Item(149, 28).Value = 50
$ws.Cells.Item(149, 29).Value = 610
$ws.Cells.Item(150, 27).Value = 25
$ws.Cells.Item(150, 28).Value = 20
$ws.Cells.Item(150, 29).Value = 135
$ws.Cells.Item(151, 27).Value = 21
$ws.Cells.Item(151, 28).Value = 50
$ws.Cells.Item(151, 29).Value = 401
$ws.Cells.Item(151, 32).Value = 147
$ws.Cells.Item(152, 27).Value = 13
$ws.Cells.Item(152, 28).Value = 55
$ws.Cells.Item(152, 29).Value = 359
$ws.Cells.Item(152, 32).Value = 117
$ws.Cells.Item(153, 27).Value = 20
$ws.Cells.Item(153, 28).Value = 25
$ws.Cells.Item(153, 29).Value = 313
$ws.Cells.Item(153, 32).Value = 123
$ws.Cells.Item(154, 27).Value = 25
$ws.Cells.Item(154, 28).Value = 35
$ws.Cells.Item(154, 29).Value = 317
$ws.Cells.Item(154, 32).Value = 130
$ws.Cells.Item(155, 27).Value = 25
$ws.Cells.Item(155, 28).Value = 75
$ws.Cells.Item(155, 29).Value = 365
$ws.Cells.Item(156, 27).Value = 20
$ws.Cells.Item(156, 28).Value = 15
$ws.Cells.Item(156, 29).Value = 65
$ws.Cells.Item(157, 27).Value = 25
$ws.Cells.Item(157, 28).Value = 40
$ws.Cells.Item(157, 29).Value = 312
$ws.Cells.Item(157, 32).Value = 143
$ws.Cells.Item(158, 27).Value = 26
$ws.Cells.Item(158, 28).Value = 40
$ws.Cells.Item(158, 29).Value = 367
$ws.Cells.Item(158, 32).Value = 143
$ws.Cells.Item(159, 27).Value = 31
$ws.Cells.Item(159, 28).Value = 35
$ws.Cells.Item(159, 29).Value = 225
$ws.Cells.Item(160, 27).Value = 0
$ws.Cells.Item(160, 28).Value = 0
$ws.Cells.Item(160, 29).Value = 357
$ws.Cells.Item(160, 32).Value = 126
$ws.Cells.Item(161, 27).Value = 0
$ws.Cells.Item(161, 28).Value = 35
$ws.Cells.Item(161, 29).Value = 556
$ws.Cells.Item(161, 32).Value = 146
$ws.Cells.Item(162, 27).Value = 0
$ws.Cells.Item(162, 28).Value = 0
$ws.Cells.Item(162, 29).Value = 252
$ws.Cells.Item(162, 32).Value = 117
$ws.Cells.Item(163, 27).Value = 0
$ws.Cells.Item(163, 28).Value = 25
$ws.Cells.Item(163, 29).Value = 390
$ws.Cells.Item(164, 27).Value = 0
$ws.Cells.Item(164, 28).Value = 20
$ws.Cells.Item(164, 29).Value = 343
$ws.Cells.Item(164, 32).Value = 126
$ws.Cells.Item(165, 27).Value = 20
$ws.Cells.Item(165, 28).Value = 45
$ws.Cells.Item(165, 29).Value = 288
$ws.Cells.Item(166, 27).Value = 27
$ws.Cells.Item(166, 28).Value = 0
$ws.Cells.Item(166, 29).Value = 246
$ws.Cells.Item(167, 27).Value = 20
$ws.Cells.Item(167, 28).Value = 0
$ws.Cells.Item(167, 29).Value = 240
$ws.Cells.Item(168, 27).Value = 20
$ws.Cells.Item(168, 28).Value = 10
$ws.Cells.Item(168, 29).Value = 45
$ws.Cells.Item(169, 27).Value = 25
$ws.Cells.Item(169, 28).Value = 50
$ws.Cells.Item(169, 29).Value = 455
$ws.Cells.Item(170, 27).Value = 20
$ws.Cells.Item(170, 28).Value = 50
$ws.Cells.Item(170, 29).Value = 383
$ws.Cells.Item(170, 32).Value = 138
$ws.Cells.Item(171, 27).Value = 0
$ws.Cells.Item(171, 28).Value = 5
$ws.Cells.Item(171, 29).Value = 253
$ws.Cells.Item(172, 27).Value = 25.5
$ws.Cells.Item(172, 28).Value = 25
$ws.Cells.Item(172, 29).Value = 120
$ws.Cells.Item(173, 27).Value = 25
$ws.Cells.Item(173, 28).Value = 10
$ws.Cells.Item(173, 29).Value = 105
$ws.Cells.Item(174, 27).Value = 20
$ws.Cells.Item(174, 28).Value = 5
$ws.Cells.Item(174, 29).Value = 30
$ws.Cells.Item(175, 27).Value = 20
$ws.Cells.Item(175, 28).Value = 35
$ws.Cells.Item(175, 29).Value = 372
$ws.Cells.Item(175, 32).Value = 137
$ws.Cells.Item(176, 27).Value = 20
$ws.Cells.Item(176, 28).Value = 45
$ws.Cells.Item(176, 29).Value = 413
$ws.Cells.Item(176, 32).Value = 144
$ws.Cells.Item(177, 27).Value = 20
$ws.Cells.Item(177, 28).Value = 15
$ws.Cells.Item(177, 29).Value = 60
$ws.Cells.Item(178, 27).Value = 20
$ws.Cells.Item(178, 28).Value = 5
$ws.Cells.Item(178, 29).Value = 40
$ws.Cells.Item(179, 27).Value = 20.5
$ws.Cells.Item(179, 28).Value = 40
$ws.Cells.Item(179, 29).Value = 283
$ws.Cells.Item(179, 32).Value = 130
$ws.Cells.Item(180, 27).Value = 20
$ws.Cells.Item(180, 28).Value = 50
$ws.Cells.Item(180, 29).Value = 487
$ws.Cells.Item(180, 32).Value = 147
$ws.Cells.Item(181, 27).Value = 20
$ws.Cells.Item(181, 28).Value = 20
$ws.Cells.Item(181, 29).Value = 80
$ws.Cells.Item(182, 27).Value = 20
$ws.Cells.Item(182, 28).Value = 35
$ws.Cells.Item(182, 29).Value = 233
$ws.Cells.Item(182, 32).Value = 126
$ws.Cells.Item(183, 27).Value = 27
$ws.Cells.Item(183, 28).Value = 40
$ws.Cells.Item(183, 29).Value = 510
$ws.Cells.Item(184, 27).Value = 15
$ws.Cells.Item(184, 28).Value = 35
$ws.Cells.Item(184, 29).Value = 273
$ws.Cells.Item(185, 27).Value = 20
$ws.Cells.Item(185, 28).Value = 25
$ws.Cells.Item(185, 29).Value = 120
$ws.Cells.Item(186, 27).Value = 22
$ws.Cells.Item(186, 28).Value = 37
$ws.Cells.Item(186, 29).Value = 377
$ws.Cells.Item(186, 32).Value = 141
$ws.Cells.Item(187, 27).Value = 25.5
$ws.Cells.Item(187, 28).Value = 25
$ws.Cells.Item(187, 29).Value = 276
$ws.Cells.Item(187, 32).Value = 134
$ws.Cells.Item(188, 27).Value = 22
$ws.Cells.Item(188, 28).Value = 40
$ws.Cells.Item(188, 29).Value = 413
$ws.Cells.Item(188, 32).Value = 137

# --- Row 189: clear all data except AF189 (CRI-Total) ---
$ws.Range("A189:Z189").ClearContents()
$ws.Range("AG189:AJ189").ClearContents()

# --- Restore selection/view state ---
$ws.Range("A1:AJ189").Select() | Out-Null
